# The sheet currently has an unused (empty) row 1, with the header row
# ("report_date" / "rate_per_100000") sitting on row 2 and the daily
# date/rate data filling rows 3-427.
#
# Delete row 1: Excel shifts every row up by one, so the header lands on
# row 1 and all the data rows move up to rows 2-426 (values/styles travel
# with their rows, only the row numbers/cell refs change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Rows("1").Delete()

# The view had scrolled down (topLeftCell=A420) with E429 selected; after
# the edit the view is back at the top with E8 selected.
[void]$ws.Range("E8").Select()
